$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(4)

$oldText = "The right side of the rule is a sequence of terminal symbols, nonterminal symbols, and other special symbols as define on the next slide."
$newText = "The right side of the rule is a sequence of terminal symbols, nonterminal symbols, and other special symbols as defined on the next slide."

if ($para.Text.TrimEnd("`r") -eq $oldText) {
    # Fix the typo ("define" -> "defined") first, while the paragraph is
    # still a single run, so the rest of the sentence is untouched.
    $run1 = $para.Runs(1)
    $run1.Text = $newText

    # Now split the corrected sentence into three runs so that the phrase
    # "as defined " (the part that actually changed) becomes its own run,
    # separate from the unedited text before and after it.
    $full = $para.Text
    $middleStart = $full.IndexOf("as defined ") + 1
    $middleLen = "as defined ".Length

    $middle = $para.Characters($middleStart, $middleLen)
    $middle.Text = "as defined "
}
